$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 729; this shifts the existing rows 729-770 down to 730-771
# (and all of their content/formatting) automatically.
$ws.Rows.Item(729).Insert()

# Populate the new row 729 with the inserted record: 2026/01/26, 月, 19, 161
# Columns A/B hold text (date-as-text / weekday kanji) in this sheet, so force
# text number-format before assigning to avoid Excel's automatic date parsing,
# then restore the default "Normal" style/format so the new cells match the
# unstyled look of the surrounding data rows.
$ws.Range("A729:B729").NumberFormat = "@"
$ws.Range("A729").Value = "2026/01/26"
$ws.Range("B729").Value = "月"
$ws.Range("A729:B729").NumberFormat = "General"
$ws.Range("A729:B729").Style = "Normal"

$ws.Range("C729").Value = 19
$ws.Range("D729").Value = 161
